# ---------------------------------------------------------------------------
# Edit summary
# ---------------------------------------------------------------------------
# 1) Slide 16 contains a table (graphicFrame) whose table style was changed
#    to a different built-in PowerPoint table style
#    ({16856D94-28A8-4DE0-957B-36C0C27DD17C}).
# 2) The deck's theme colour scheme was changed from the "Integral" palette
#    to the standard "Office" palette (this is what drives
#    ppt/theme/theme1.xml, the theme used by the slide master / all slides).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Apply the new table style to the table on slide 16.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(16)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{16856D94-28A8-4DE0-957B-36C0C27DD17C}")
    }
}

# ---------------------------------------------------------------------------
# 2) Re-colour the presentation's theme colour scheme to the standard
#    "Office" palette. The colour scheme is shared by the whole deck, so it
#    can be reached from any slide's ThemeColorScheme.
#
#    ThemeColorScheme.Colors index order:
#      1 dk1   2 lt1   3 dk2   4 lt2
#      5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#      11 hlink   12 folHlink
#
#    .RGB uses the classic VBA BGR-packed integer (R + G*256 + B*65536),
#    i.e. the same encoding produced by VBA's RGB(r,g,b) function.
# ---------------------------------------------------------------------------
$themeColors = $p.Slides.Item(1).ThemeColorScheme

# dk1  -> 000000
$themeColors.Colors(1).RGB = 0
# lt1  -> FFFFFF
$themeColors.Colors(2).RGB = 16777215
# dk2  -> 44546A
$themeColors.Colors(3).RGB = 6968388
# lt2  -> E7E6E6
$themeColors.Colors(4).RGB = 15132391
# accent1 -> 5B9BD5
$themeColors.Colors(5).RGB = 13998939
# accent2 -> ED7D31
$themeColors.Colors(6).RGB = 3243501
# accent3 -> A5A5A5
$themeColors.Colors(7).RGB = 10855845
# accent4 -> FFC000
$themeColors.Colors(8).RGB = 49407
# accent5 -> 4472C4
$themeColors.Colors(9).RGB = 12874308
# accent6 -> 70AD47
$themeColors.Colors(10).RGB = 4697456
# hlink -> 0563C1
$themeColors.Colors(11).RGB = 12673797
# folHlink -> 954F72
$themeColors.Colors(12).RGB = 7491477
